# Auto-generated Excel COM-interop script to apply Raiden_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1007.0294  # H17
$ws.Cells.Item(17, 10).Value = 1007.0294  # J17
$ws.Cells.Item(17, 12).Value = 3021.0882  # L17
$ws.Cells.Item(17, 14).Value = -3357.0882  # N17
$ws.Cells.Item(28, 8).Value = 1330.8422  # H28
$ws.Cells.Item(28, 9).Value = 1099.5714  # I28
$ws.Cells.Item(28, 11).Value = 1099.5714  # K28
$ws.Cells.Item(28, 13).Value = -614.5714  # M28
$ws.Cells.Item(30, 8).Value = 0  # H30
$ws.Cells.Item(30, 10).Value = 0  # J30
$ws.Cells.Item(30, 12).Value = 0  # L30
$ws.Cells.Item(30, 14).ClearContents() | Out-Null  # N30
$ws.Cells.Item(40, 8).Value = 2675.1  # H40
$ws.Cells.Item(40, 10).Value = 2858  # J40
$ws.Cells.Item(40, 12).Value = 2858  # L40
$ws.Cells.Item(40, 14).Value = -3208  # N40
$ws.Cells.Item(58, 8).Value = 2098.5  # H58
$ws.Cells.Item(58, 10).Value = 999  # J58
$ws.Cells.Item(58, 12).Value = 2997  # L58
$ws.Cells.Item(58, 14).Value = -3297  # N58
$ws.Cells.Item(76, 8).Value = 4863  # H76
$ws.Cells.Item(76, 9).Value = 3898.2  # I76
$ws.Cells.Item(76, 10).Value = 5399  # J76
$ws.Cells.Item(76, 11).Value = 3898.2  # K76
$ws.Cells.Item(76, 12).Value = 5399  # L76
$ws.Cells.Item(76, 13).Value = -3583.2  # M76
$ws.Cells.Item(76, 14).Value = -6029  # N76
$ws.Cells.Item(79, 8).Value = 4863  # H79
$ws.Cells.Item(79, 9).Value = 3898.2  # I79
$ws.Cells.Item(79, 10).Value = 5399  # J79
$ws.Cells.Item(79, 11).Value = 3898.2  # K79
$ws.Cells.Item(79, 12).Value = 5399  # L79
$ws.Cells.Item(79, 13).Value = -2806.2  # M79
$ws.Cells.Item(79, 14).Value = -7583  # N79
$ws.Cells.Item(106, 8).Value = 1637  # H106
$ws.Cells.Item(106, 9).Value = 1442.2858  # I106
$ws.Cells.Item(106, 11).Value = 1442.2858  # K106
$ws.Cells.Item(106, 13).Value = -811.2858000000001  # M106
$ws.Cells.Item(107, 8).Value = 1045.5385  # H107
$ws.Cells.Item(107, 9).Value = 1129  # I107
$ws.Cells.Item(107, 11).Value = 1129  # K107
$ws.Cells.Item(107, 13).Value = 791  # M107
$ws.Cells.Item(132, 8).Value = 910997  # H132
$ws.Cells.Item(132, 9).Value = 2107.4443  # I132
$ws.Cells.Item(132, 10).Value = 5001000  # J132
$ws.Cells.Item(132, 11).Value = 6322.3329  # K132
$ws.Cells.Item(132, 12).Value = 15003000  # L132
$ws.Cells.Item(132, 13).Value = -3792.3329  # M132
$ws.Cells.Item(132, 14).Value = -15008060  # N132
$ws.Cells.Item(135, 8).Value = 5271.1304  # H135
$ws.Cells.Item(135, 9).Value = 936.85  # I135
$ws.Cells.Item(135, 11).Value = 8431.65  # K135
$ws.Cells.Item(135, 13).Value = -5896.65  # M135
$ws.Cells.Item(138, 8).Value = 2599.19  # H138
$ws.Cells.Item(138, 10).Value = 2146.5  # J138
$ws.Cells.Item(138, 12).Value = 6439.5  # L138
$ws.Cells.Item(138, 14).Value = -16719.5  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2590.9678  # H61
$ws.Cells.Item(61, 9).Value = 2512.6538  # I61
$ws.Cells.Item(61, 11).Value = 2512.6538  # K61
$ws.Cells.Item(61, 13).Value = -2300.6538  # M61
$ws.Cells.Item(74, 8).Value = 1687.15  # H74
$ws.Cells.Item(74, 9).Value = 1652.5555  # I74
$ws.Cells.Item(74, 11).Value = 1652.5555  # K74
$ws.Cells.Item(74, 13).Value = -778.5554999999999  # M74
$ws.Cells.Item(77, 8).Value = 1687.15  # H77
$ws.Cells.Item(77, 9).Value = 1652.5555  # I77
$ws.Cells.Item(77, 11).Value = 8262.7775  # K77
$ws.Cells.Item(77, 13).Value = -3894.7775  # M77
$ws.Cells.Item(97, 8).Value = 997.6111  # H97
$ws.Cells.Item(97, 9).Value = 325.13333  # I97
$ws.Cells.Item(97, 10).Value = 4360  # J97
$ws.Cells.Item(97, 11).Value = 325.13333  # K97
$ws.Cells.Item(97, 12).Value = 4360  # L97
$ws.Cells.Item(97, 13).Value = 170.86667  # M97
$ws.Cells.Item(97, 14).Value = -5352  # N97
$ws.Cells.Item(132, 8).Value = 2528.0625  # H132
$ws.Cells.Item(132, 9).Value = 2004.0834  # I132
$ws.Cells.Item(132, 10).Value = 4100  # J132
$ws.Cells.Item(132, 11).Value = 6012.2502  # K132
$ws.Cells.Item(132, 12).Value = 12300  # L132
$ws.Cells.Item(132, 13).Value = -3482.2502  # M132
$ws.Cells.Item(132, 14).Value = -17360  # N132
$ws.Cells.Item(136, 8).Value = 2590.9678  # H136
$ws.Cells.Item(136, 9).Value = 2512.6538  # I136
$ws.Cells.Item(136, 11).Value = 7537.9614  # K136
$ws.Cells.Item(136, 13).Value = -4987.9614  # M136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(63, 8).Value = 37635.5  # H63
$ws.Cells.Item(63, 9).Value = 5000  # I63
$ws.Cells.Item(63, 10).Value = 70271  # J63
$ws.Cells.Item(63, 11).Value = 5000  # K63
$ws.Cells.Item(63, 12).Value = 70271  # L63
$ws.Cells.Item(63, 13).Value = -4314  # M63
$ws.Cells.Item(63, 14).Value = -71643  # N63
$ws.Cells.Item(66, 8).Value = 37635.5  # H66
$ws.Cells.Item(66, 9).Value = 5000  # I66
$ws.Cells.Item(66, 10).Value = 70271  # J66
$ws.Cells.Item(66, 11).Value = 15000  # K66
$ws.Cells.Item(66, 12).Value = 210813  # L66
$ws.Cells.Item(66, 13).Value = -11568  # M66
$ws.Cells.Item(66, 14).Value = -217677  # N66
$ws.Cells.Item(96, 8).Value = 18724.5  # H96
$ws.Cells.Item(96, 9).Value = 14469.4  # I96
$ws.Cells.Item(96, 11).Value = 14469.4  # K96
$ws.Cells.Item(96, 13).Value = -11723.4  # M96
$ws.Cells.Item(132, 8).Value = 0  # H132
$ws.Cells.Item(132, 10).Value = 0  # J132
$ws.Cells.Item(132, 12).Value = 0  # L132
$ws.Cells.Item(132, 14).ClearContents() | Out-Null  # N132

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2340.48  # H31
$ws.Cells.Item(31, 9).Value = 1671  # I31
$ws.Cells.Item(31, 11).Value = 1671  # K31
$ws.Cells.Item(31, 13).Value = -1376  # M31
$ws.Cells.Item(34, 8).Value = 2340.48  # H34
$ws.Cells.Item(34, 9).Value = 1671  # I34
$ws.Cells.Item(34, 11).Value = 1671  # K34
$ws.Cells.Item(34, 13).Value = -1469  # M34
$ws.Cells.Item(41, 8).Value = 5000  # H41
$ws.Cells.Item(41, 10).Value = 0  # J41
$ws.Cells.Item(41, 12).Value = 0  # L41
$ws.Cells.Item(41, 14).ClearContents() | Out-Null  # N41
$ws.Cells.Item(60, 8).Value = 22082.334  # H60
$ws.Cells.Item(60, 10).Value = 29998  # J60
$ws.Cells.Item(60, 12).Value = 29998  # L60
$ws.Cells.Item(60, 14).Value = -31020  # N60
$ws.Cells.Item(107, 8).Value = 1775.8889  # H107
$ws.Cells.Item(107, 9).Value = 1747.875  # I107
$ws.Cells.Item(107, 11).Value = 1747.875  # K107
$ws.Cells.Item(107, 13).Value = 172.125  # M107

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 500055.6  # H2
$ws.Cells.Item(2, 9).Value = 625048.1  # I2
$ws.Cells.Item(2, 11).Value = 3750288.6  # K2
$ws.Cells.Item(2, 13).Value = -3750175.6  # M2
$ws.Cells.Item(4, 8).Value = 57116444  # H4
$ws.Cells.Item(4, 9).Value = 3864787.8  # I4
$ws.Cells.Item(4, 10).Value = 350000540  # J4
$ws.Cells.Item(4, 11).Value = 11594363.4  # K4
$ws.Cells.Item(4, 12).Value = 1050001620  # L4
$ws.Cells.Item(4, 13).Value = -11594251.4  # M4
$ws.Cells.Item(4, 14).Value = -1050001844  # N4
$ws.Cells.Item(12, 8).Value = 267.5  # H12
$ws.Cells.Item(12, 9).Value = 327.5  # I12
$ws.Cells.Item(12, 10).Value = 216.07143  # J12
$ws.Cells.Item(12, 11).Value = 982.5  # K12
$ws.Cells.Item(12, 12).Value = 648.21429  # L12
$ws.Cells.Item(12, 13).Value = -809.5  # M12
$ws.Cells.Item(12, 14).Value = -994.21429  # N12
$ws.Cells.Item(14, 8).Value = 201.28572  # H14
$ws.Cells.Item(14, 9).Value = 201.28572  # I14
$ws.Cells.Item(14, 11).Value = 603.85716  # K14
$ws.Cells.Item(14, 13).Value = -430.85716  # M14
$ws.Cells.Item(97, 8).Value = 1084.125  # H97
$ws.Cells.Item(97, 10).Value = 574.75  # J97
$ws.Cells.Item(97, 12).Value = 1724.25  # L97
$ws.Cells.Item(97, 14).Value = -2716.25  # N97
$ws.Cells.Item(114, 8).Value = 1640.2858  # H114
$ws.Cells.Item(114, 9).Value = 1526.2858  # I114
$ws.Cells.Item(114, 10).Value = 1754.2858  # J114
$ws.Cells.Item(114, 11).Value = 4578.857400000001  # K114
$ws.Cells.Item(114, 12).Value = 5262.857400000001  # L114
$ws.Cells.Item(114, 13).Value = -1324.857400000001  # M114
$ws.Cells.Item(114, 14).Value = -11770.8574  # N114
$ws.Cells.Item(132, 8).Value = 952.8182  # H132
$ws.Cells.Item(132, 9).Value = 747.625  # I132
$ws.Cells.Item(132, 11).Value = 6728.625  # K132
$ws.Cells.Item(132, 13).Value = -4198.625  # M132

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 591.25  # H107
$ws.Cells.Item(107, 9).Value = 241.11111  # I107
$ws.Cells.Item(107, 11).Value = 241.11111  # K107
$ws.Cells.Item(107, 13).Value = 1678.88889  # M107

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3717.7273  # H46
$ws.Cells.Item(46, 9).Value = 2448  # I46
$ws.Cells.Item(46, 10).Value = 3999.889  # J46
$ws.Cells.Item(46, 11).Value = 2448  # K46
$ws.Cells.Item(46, 12).Value = 3999.889  # L46
$ws.Cells.Item(46, 13).Value = -2260  # M46
$ws.Cells.Item(46, 14).Value = -4375.889  # N46
$ws.Cells.Item(82, 8).Value = 2893.8333  # H82
$ws.Cells.Item(82, 9).Value = 2341  # I82
$ws.Cells.Item(82, 10).Value = 3999.5  # J82
$ws.Cells.Item(82, 11).Value = 2341  # K82
$ws.Cells.Item(82, 12).Value = 3999.5  # L82
$ws.Cells.Item(82, 13).Value = -1980  # M82
$ws.Cells.Item(82, 14).Value = -4721.5  # N82
$ws.Cells.Item(85, 8).Value = 2893.8333  # H85
$ws.Cells.Item(85, 9).Value = 2341  # I85
$ws.Cells.Item(85, 10).Value = 3999.5  # J85
$ws.Cells.Item(85, 11).Value = 2341  # K85
$ws.Cells.Item(85, 12).Value = 3999.5  # L85
$ws.Cells.Item(85, 13).Value = -1093  # M85
$ws.Cells.Item(85, 14).Value = -6495.5  # N85

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(68, 8).Value = 12122.5  # H68
$ws.Cells.Item(68, 9).Value = 14245  # I68
$ws.Cells.Item(68, 10).Value = 10000  # J68
$ws.Cells.Item(68, 11).Value = 14245  # K68
$ws.Cells.Item(68, 12).Value = 10000  # L68
$ws.Cells.Item(68, 13).Value = -13434  # M68
$ws.Cells.Item(68, 14).Value = -11622  # N68
$ws.Cells.Item(69, 8).Value = 36138  # H69
$ws.Cells.Item(69, 10).Value = 36138  # J69
$ws.Cells.Item(69, 12).Value = 36138  # L69
$ws.Cells.Item(69, 14).Value = -37636  # N69
$ws.Cells.Item(70, 8).Value = 18674.65  # H70
$ws.Cells.Item(70, 9).Value = 14749.167  # I70
$ws.Cells.Item(70, 10).Value = 20357  # J70
$ws.Cells.Item(70, 11).Value = 14749.167  # K70
$ws.Cells.Item(70, 12).Value = 20357  # L70
$ws.Cells.Item(70, 13).Value = -14434.167  # M70
$ws.Cells.Item(70, 14).Value = -20987  # N70
$ws.Cells.Item(71, 8).Value = 12122.5  # H71
$ws.Cells.Item(71, 9).Value = 14245  # I71
$ws.Cells.Item(71, 10).Value = 10000  # J71
$ws.Cells.Item(71, 11).Value = 42735  # K71
$ws.Cells.Item(71, 12).Value = 30000  # L71
$ws.Cells.Item(71, 13).Value = -38679  # M71
$ws.Cells.Item(71, 14).Value = -38112  # N71
$ws.Cells.Item(72, 8).Value = 36138  # H72
$ws.Cells.Item(72, 10).Value = 36138  # J72
$ws.Cells.Item(72, 12).Value = 108414  # L72
$ws.Cells.Item(72, 14).Value = -115902  # N72
$ws.Cells.Item(73, 8).Value = 18674.65  # H73
$ws.Cells.Item(73, 9).Value = 14749.167  # I73
$ws.Cells.Item(73, 10).Value = 20357  # J73
$ws.Cells.Item(73, 11).Value = 14749.167  # K73
$ws.Cells.Item(73, 12).Value = 20357  # L73
$ws.Cells.Item(73, 13).Value = -13657.167  # M73
$ws.Cells.Item(73, 14).Value = -22541  # N73
$ws.Cells.Item(99, 8).Value = 35249.25  # H99
$ws.Cells.Item(99, 9).Value = 30332.334  # I99
$ws.Cells.Item(99, 11).Value = 30332.334  # K99
$ws.Cells.Item(99, 13).Value = -27337.334  # M99

"Done: 228 sets, 3 clears"